# Generate Report for Handoff
# Updates the "aeb9f62d-7abc-4fe6-81aa-dc8447d2da20.md" row on all three sheets
# to reflect the file having been handed off for translation.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Overview sheet: row 3 is the aeb9f62d... file.
#   E (zh-cn status) / F (de-de status): "In Translation" -> "Ready for handoff"
#   G (Latest HO Xliff Generate Date): "2016-08-20 22:12:38" -> "2016-08-20 22:13:16"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-20 22:13:16"

# ---------------------------------------------------------------------------
# zh-cn sheet: row 3 is the aeb9f62d... file.
#   C (Status): "In Translation" -> "Ready for handoff"
#   E (Priority): "ht" -> "mt"
#   H (Latest Handoff Datetime): "2016-08-20 22:12:34" -> "2016-08-20 22:13:12"
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("H3").Value = "2016-08-20 22:13:12"

# ---------------------------------------------------------------------------
# de-de sheet: row 3 is the aeb9f62d... file.
#   C (Status): "In Translation" -> "Ready for handoff"
#   E (Priority): "ht" -> "mt"
#   H (Latest Handoff Datetime): "2016-08-20 22:12:38" -> "2016-08-20 22:13:16"
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("H3").Value = "2016-08-20 22:13:16"

# ---------------------------------------------------------------------------
# The longer "Ready for handoff" status text (vs. "In Translation") makes
# Excel auto-fit the Status/zh-cn/de-de columns a bit wider on every sheet.
# ---------------------------------------------------------------------------
$wsOverview.Columns.Item(5).AutoFit()
$wsOverview.Columns.Item(6).AutoFit()
$wsZhCn.Columns.Item(3).AutoFit()
$wsDeDe.Columns.Item(3).AutoFit()
